$wb = $excel.ActiveWorkbook

# --- Sheet "Senador": add a couple of new grade entries for two students ---
$senador = $wb.Worksheets.Item("Senador")

# Juan Pablo Alfaya (row 3): new scores of 0 in columns D and E
$senador.Range("D3").Value = 0
$senador.Range("E3").Value = 0

# Row 6 student: new score of 0 in column D
$senador.Range("D6").Value = 0

# Leave the active selection where the edits were last made
$senador.Range("D7").Select()

# --- Sheet "Mago": add a new grade entry for the student in row 4 ---
$mago = $wb.Worksheets.Item("Mago")

# Row 4 student: new score of 1 in column H
$mago.Range("H4").Value = 1

# Leave the active selection where the edit was last made
$mago.Range("H5").Select()
